$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 87.083336
$ws.Range("I2").Value = 87.083336
$ws.Range("K2").Value = 87.083336
$ws.Range("M2").Value = 25.916664
$ws.Range("H4").Value = 129.90909
$ws.Range("I4").Value = 102.71429
$ws.Range("J4").Value = 177.5
$ws.Range("K4").Value = 102.71429
$ws.Range("L4").Value = 177.5
$ws.Range("M4").Value = 11.28570999999999
$ws.Range("N4").Value = -405.5
$ws.Range("H29").Value = 725
$ws.Range("I29").Value = 587.5
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 1762.5
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = -1481.5
$ws.Range("N29").Value = -3562
$ws.Range("H32").Value = 1546.3334
$ws.Range("J32").Value = 1455.6
$ws.Range("L32").Value = 1455.6
$ws.Range("N32").Value = -2107.6
$ws.Range("H38").Value = 429.8
$ws.Range("I38").Value = 62.25
$ws.Range("J38").Value = 1900
$ws.Range("K38").Value = 186.75
$ws.Range("L38").Value = 5700
$ws.Range("M38").Value = 185.25
$ws.Range("N38").Value = -6444
$ws.Range("H39").Value = 625.8
$ws.Range("I39").Value = 532.25
$ws.Range("K39").Value = 1596.75
$ws.Range("M39").Value = -1300.75
$ws.Range("H43").Value = 3498.8572
$ws.Range("I43").Value = 3498.8572
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3498.8572
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3429.8572
$ws.Range("N43").Value = ""
$ws.Range("H58").Value = 308.33334
$ws.Range("I58").Value = 308.33334
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 925.0000200000001
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -775.0000200000001
$ws.Range("N58").Value = ""
$ws.Range("H70").Value = 1310.5555
$ws.Range("J70").Value = 1549
$ws.Range("L70").Value = 4647
$ws.Range("N70").Value = -5187
$ws.Range("H73").Value = 1310.5555
$ws.Range("J73").Value = 1549
$ws.Range("L73").Value = 4647
$ws.Range("N73").Value = -6519
$ws.Range("H86").Value = 4310
$ws.Range("I86").Value = 2997.7
$ws.Range("J86").Value = 5768.1113
$ws.Range("K86").Value = 2997.7
$ws.Range("L86").Value = 5768.1113
$ws.Range("M86").Value = -1874.7
$ws.Range("N86").Value = -8014.1113
$ws.Range("H89").Value = 4310
$ws.Range("I89").Value = 2997.7
$ws.Range("J89").Value = 5768.1113
$ws.Range("K89").Value = 14988.5
$ws.Range("L89").Value = 28840.5565
$ws.Range("M89").Value = -9372.5
$ws.Range("N89").Value = -40072.5565
$ws.Range("H98").Value = 589
$ws.Range("I98").Value = 589
$ws.Range("K98").Value = 589
$ws.Range("M98").Value = 909
$ws.Range("H122").Value = 589
$ws.Range("I122").Value = 589
$ws.Range("K122").Value = 1767
$ws.Range("M122").Value = 683
$ws.Range("H137").Value = 1262
$ws.Range("I137").Value = 1386
$ws.Range("K137").Value = 4158
$ws.Range("M137").Value = -1608
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1272.75
$ws.Range("I5").Value = 14
$ws.Range("K5").Value = 14
$ws.Range("M5").Value = 98
$ws.Range("H12").Value = 6147.5
$ws.Range("I12").Value = 2300
$ws.Range("J12").Value = 9995
$ws.Range("K12").Value = 2300
$ws.Range("L12").Value = 9995
$ws.Range("M12").Value = -2127
$ws.Range("N12").Value = -10341
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1272.75
$ws.Range("I4").Value = 14
$ws.Range("K4").Value = 14
$ws.Range("M4").Value = 101
$ws.Range("H105").Value = 2151.1667
$ws.Range("I105").Value = 2151.1667
$ws.Range("K105").Value = 2151.1667
$ws.Range("M105").Value = -404.1667000000002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 30.666666
$ws.Range("I6").Value = 90
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 90
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 23
$ws.Range("N6").Value = -227
$ws.Range("H58").Value = 676.5833
$ws.Range("I58").Value = 511.9
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 511.9
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -308.9
$ws.Range("N58").Value = -1906
$ws.Range("H86").Value = 402266.34
$ws.Range("I86").Value = 242719.8
$ws.Range("K86").Value = 242719.8
$ws.Range("M86").Value = -241596.8
$ws.Range("H89").Value = 402266.34
$ws.Range("I89").Value = 242719.8
$ws.Range("K89").Value = 1213599
$ws.Range("M89").Value = -1207983
$ws.Range("H94").Value = 3327.6
$ws.Range("J94").Value = 2986
$ws.Range("L94").Value = 2986
$ws.Range("N94").Value = -3888
$ws.Range("H99").Value = 3945.2942
$ws.Range("I99").Value = 4326.4287
$ws.Range("J99").Value = 2166.6667
$ws.Range("K99").Value = 4326.4287
$ws.Range("L99").Value = 2166.6667
$ws.Range("M99").Value = -2828.4287
$ws.Range("N99").Value = -5162.6667
$ws.Range("H126").Value = 3945.2942
$ws.Range("I126").Value = 4326.4287
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 12979.2861
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -10509.2861
$ws.Range("N126").Value = -11440.0001
$ws.Range("H132").Value = 3415.5454
$ws.Range("J132").Value = 1610
$ws.Range("L132").Value = 4830
$ws.Range("N132").Value = -9890
$ws.Range("H134").Value = 1857.9166
$ws.Range("I134").Value = 1726.7
$ws.Range("K134").Value = 5180.1
$ws.Range("M134").Value = -2645.1
$ws.Range("H136").Value = 676.5833
$ws.Range("I136").Value = 511.9
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 1535.7
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = 1014.3
$ws.Range("N136").Value = -9600
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 105.5
$ws.Range("I6").Value = 105.5
$ws.Range("K6").Value = 316.5
$ws.Range("M6").Value = -203.5
$ws.Range("H12").Value = 1268.75
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1268.75
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 3806.25
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = -4152.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10002
$ws.Range("I5").Value = 10000
$ws.Range("J5").Value = 10004
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10004
$ws.Range("M5").Value = -9888
$ws.Range("N5").Value = -10228
$ws.Range("H113").Value = 10011
$ws.Range("I113").Value = 10011
$ws.Range("K113").Value = 10011
$ws.Range("M113").Value = -7841
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -888
$ws.Range("N2").Value = ""
$ws.Range("H9").Value = 368
$ws.Range("I9").Value = 431.4
$ws.Range("J9").Value = 51
$ws.Range("K9").Value = 431.4
$ws.Range("L9").Value = 51
$ws.Range("M9").Value = -207.4
$ws.Range("N9").Value = -499
$ws.Range("H20").Value = 1255.5
$ws.Range("I20").Value = 505
$ws.Range("J20").Value = 2006
$ws.Range("K20").Value = 505
$ws.Range("L20").Value = 2006
$ws.Range("M20").Value = -279
$ws.Range("N20").Value = -2458
$ws.Range("H35").Value = 1399.6666
$ws.Range("I35").Value = 1399.6666
$ws.Range("K35").Value = 1399.6666
$ws.Range("M35").Value = -1063.6666
$ws.Range("H132").Value = 1857.1538
$ws.Range("I132").Value = 1844.8334
$ws.Range("K132").Value = 5534.5002
$ws.Range("M132").Value = -3004.5002
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 750000
$ws.Range("I3").Value = 1000000
$ws.Range("J3").Value = 500000
$ws.Range("K3").Value = 1000000
$ws.Range("L3").Value = 500000
$ws.Range("M3").Value = -999886
$ws.Range("N3").Value = -500228
$ws.Range("H5").Value = 600000000
$ws.Range("I5").Value = 600000000
$ws.Range("K5").Value = 600000000
$ws.Range("M5").Value = -599999888
$ws.Range("H62").Value = 1499.5
$ws.Range("J62").Value = 1499.5
$ws.Range("L62").Value = 1499.5
$ws.Range("N62").Value = -2747.5
$ws.Range("H65").Value = 1499.5
$ws.Range("J65").Value = 1499.5
$ws.Range("L65").Value = 7497.5
$ws.Range("N65").Value = -13737.5
$ws.Range("H81").Value = 5768
$ws.Range("I81").Value = 5768
$ws.Range("K81").Value = 11536
$ws.Range("M81").Value = -10475
$ws.Range("H84").Value = 5768
$ws.Range("I84").Value = 5768
$ws.Range("K84").Value = 57680
$ws.Range("M84").Value = -52376
$ws.Range("H119").Value = 48399.5
$ws.Range("J119").Value = 48399.5
$ws.Range("L119").Value = 48399.5
$ws.Range("N119").Value = -58075.5
$ws.Range("H141").Value = 71250
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 71250
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 71250
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -81610
